# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 22:22"

# Row 32 / 33 - Tenerife and Burgos swap order/data
$ws.Range("A32").Value = "Burgos"
$ws.Range("B32").Value = 1232
$ws.Range("C32").Value = 517
$ws.Range("D32").Value = 574
$ws.Range("E32").Value = 141

$ws.Range("A33").Value = "Tenerife"
$ws.Range("B33").Value = 1220
$ws.Range("C33").Value = 311
$ws.Range("D33").Value = 840
$ws.Range("E33").Value = 69

# Row 49 - Gran Canaria data refresh
$ws.Range("B49").Value = 455
$ws.Range("C49").Value = 189
$ws.Range("D49").Value = 239
$ws.Range("E49").Value = 27

# Row 55 - Ceuta data refresh
$ws.Range("B55").Value = 96

# Row 56 - La Palma data refresh
$ws.Range("B56").Value = 83
$ws.Range("C56").Value = 14
$ws.Range("D56").Value = 66

# Row 57 - Lanzarote data refresh
$ws.Range("B57").Value = 79
$ws.Range("C57").Value = 19
$ws.Range("D57").Value = 57
$ws.Range("E57").Value = 3

# Row 59 - Fuerteventura data refresh
$ws.Range("B59").Value = 42
$ws.Range("D59").Value = 34

# Row 62 / 63 - Arroyo de la Luz and La Gomera swap order/data
$ws.Range("A62").Value = "La Gomera"
$ws.Range("B62").Value = 10
$ws.Range("C62").Value = 7
$ws.Range("D62").Value = 3

$ws.Range("A63").Value = "Arroyo de la Luz"
$ws.Range("C63").Value = 0
$ws.Range("D63").Value = 7

# Row 64 - El Hierro data refresh
$ws.Range("B64").Value = 4
$ws.Range("C64").Value = 2
$ws.Range("D64").Value = 2
